# Commit: "Change to MethylAction case"
# Rebrands "Methylaction" -> "MethylAction" in a handful of prose runs,
# bumps the document date, and adds quotes around the literal R package
# name ("methylaction") in one sentence. The many other lower-case
# "methylaction" occurrences (package/function/file/folder names such as
# methylaction(), methylaction_demo.tar.gz, jeffbhasin/methylaction, the
# library(methylaction) call, etc.) are intentionally left untouched.

$d = $word.ActiveDocument

# Turn off smart-quote autocorrect so that any literal " we insert below
# stays a straight quote instead of being turned into curly quotes.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# 1. Title: "Methylaction:" -> "MethylAction:"
$r1 = $d.Content.Find.Execute("Methylaction:", $true, $true, $false, $false, $false, $true, 1, $false, "MethylAction:", 2)
Write-Host "1 title: $r1"

# 2. Date: "2015-05-07" -> "2015-05-25"
$r2 = $d.Content.Find.Execute("2015-05-07", $true, $true, $false, $false, $false, $true, 1, $false, "2015-05-25", 2)
Write-Host "2 date: $r2"

# 3. Purpose paragraph has two occurrences of "Methylaction" -> "MethylAction"
$r3a = $d.Content.Find.Execute("the use of Methylaction to detect", $true, $true, $false, $false, $false, $true, 1, $false, "the use of MethylAction to detect", 2)
Write-Host "3a purpose (use of): $r3a"
$r3b = $d.Content.Find.Execute("While Methylaction is designed", $true, $true, $false, $false, $false, $true, 1, $false, "While MethylAction is designed", 2)
Write-Host "3b purpose (while): $r3b"

# 4. Installation paragraph: "Goldmine and Methylaction from GitHub" -> "Goldmine and MethylAction from GitHub"
$r4 = $d.Content.Find.Execute("Goldmine and Methylaction from GitHub", $true, $true, $false, $false, $false, $true, 1, $false, "Goldmine and MethylAction from GitHub", 2)
Write-Host "4 install: $r4"

# 5. Preprocessing paragraph: add straight quotes around "methylaction"
#    (the package name keeps its lower-case spelling here, only quoting
#    is added). Locate the text with Find, then assign Range.Text
#    directly so the literal double-quote characters are written as
#    straight quotes rather than Find/Replace's curly-quote autocorrect.
$rng5 = $d.Content
$r5 = $rng5.Find.Execute("load the methylaction R package", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "5 found: $r5"
if ($r5) {
    $rng5.Text = "load the " + [char]34 + "methylaction" + [char]34 + " R package"
}

# 6. "each run of Methylaction." -> "each run of MethylAction."
$r6 = $d.Content.Find.Execute("each run of Methylaction.", $true, $true, $false, $false, $false, $true, 1, $false, "each run of MethylAction.", 2)
Write-Host "6 reads: $r6"
